$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header row: a new "id" column is inserted before the old univoc_id
# column, shifting everything else one column to the right. Column I
# is brand new, so first clone the bold/bordered header style from H1
# onto it (keeps the same shared style index instead of minting a
# duplicate one), then overwrite the value.
# ------------------------------------------------------------------
$ws.Range("H1").Copy($ws.Range("I1"))

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "univoc_id"
$ws.Range("C1").Value = "id_diary"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "description"
$ws.Range("F1").Value = "date start"
$ws.Range("G1").Value = "date and"
$ws.Range("H1").Value = "do"
$ws.Range("I1").Value = "repeat"

# ------------------------------------------------------------------
# The three "DB" rows (2-4) now store real datetimes instead of text,
# formatted with a custom date/time display format. (entered once in
# lower case, then corrected to upper case - this is why both casings
# end up registered as custom number formats, only the upper one used)
# ------------------------------------------------------------------
$ws.Range("F2:G4").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("F2:G4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ------------------------------------------------------------------
# Columns B (univoc_id) and C (id_diary) hold big numeric-looking
# identifiers that must stay text. Mark the whole block as Text before
# typing them in so Excel doesn't silently coerce them to floats/
# scientific notation, then drop the formatting again so the cells
# come out plain (matching how they were authored from the DB export).
# ------------------------------------------------------------------
$ws.Range("B2:C7").NumberFormat = "@"

$ws.Range("B2").Value = "3541791985364674716"
$ws.Range("C2").Value = "65592250285068942839"

$ws.Range("B3").Value = "6115568889517910016"
$ws.Range("C3").Value = "65592250285068942839"

$ws.Range("B4").Value = "774462689499478238"
$ws.Range("C4").Value = "65592250285068942839"

$ws.Range("B5").Value = "40764073562277591648"
$ws.Range("C5").Value = "65592250285068942839"

$ws.Range("B6").Value = "4225836343232735037"
$ws.Range("C6").Value = "65592250285068942839"

$ws.Range("B7").Value = "73313581294534868820"
$ws.Range("C7").Value = "65592250285068942839"

$ws.Range("B2:C7").ClearFormats()

# --- Row 2: Daily event ---
$ws.Range("A2").Value = 25
$ws.Range("D2").Value = "Daily event"
$ws.Range("E2").Value = "My first daily event modify"
$ws.Range("F2").Value = 45519.625
$ws.Range("G2").Value = 45519.625
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

# --- Row 3: Monthly Event ---
$ws.Range("A3").Value = 54
$ws.Range("D3").Value = "Monthly Event"
$ws.Range("E3").Value = "My first event monthly modify"
$ws.Range("F3").Value = 45474
$ws.Range("G3").Value = 45474
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0

# --- Row 4: Period event ---
$ws.Range("A4").Value = 55
$ws.Range("D4").Value = "Period event"
$ws.Range("E4").Value = "My first event period modify"
$ws.Range("F4").Value = 45139.41666666666
$ws.Range("G4").Value = 45139.41666666666
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0

# --- Row 5: Daily event2 (new row, still the old text-date / bool shape) ---
$ws.Range("D5").Value = "Daily event2"
$ws.Range("E5").Value = "My first daily event modify2"
$ws.Range("F5").Value = "2024-08-15 15:00:00"
$ws.Range("G5").Value = "2024-08-15 15:06:00"
$ws.Range("H5").Value = $true
$ws.Range("I5").Value = $false

# --- Row 6: Monthly Event2 ---
$ws.Range("D6").Value = "Monthly Event2"
$ws.Range("E6").Value = "My first event monthly modify2"
$ws.Range("F6").Value = "2024-07-01 00:00:00"
$ws.Range("G6").Value = "2024-07-31 00:00:00"
$ws.Range("H6").Value = $true
$ws.Range("I6").Value = $false

# --- Row 7: Period event2 ---
$ws.Range("D7").Value = "Period event2"
$ws.Range("E7").Value = "My first event period modify2"
$ws.Range("F7").Value = "2023-08-01 10:00:00"
$ws.Range("G7").Value = "2024-10-31 10:00:00"
$ws.Range("H7").Value = $true
$ws.Range("I7").Value = $false
